# Apply linting fixes described in commit:
#   "fixing formatting and other issues identified with linting"
#
# 1) Bump ObjTables header metadata (version 0.0.8 -> 0.0.9, date
#    2020-03-09 23:59:06 -> 2020-04-26 21:09:24) on every !!<Table> sheet's
#    title row(s).
# 2) Rename the "From"/"To" columns on the Relation / Relationship sheets to
#    "FromObject" / "ToObject" (both the header cell text and the matching
#    data-validation titles/prompts).

$wb = $excel.ActiveWorkbook

$oldDate = "date='2020-03-09 23:59:06'"
$newDate = "date='2020-04-26 21:09:24'"
$oldVer  = "objTablesVersion='0.0.8'"
$newVer  = "objTablesVersion='0.0.9'"

function Bump-ObjTables([string]$text) {
    $t = $text -replace [regex]::Escape($oldDate), $newDate
    $t = $t -replace [regex]::Escape($oldVer), $newVer
    return $t
}

# --- 1) Walk every worksheet and patch its ObjTables header string(s) ----
foreach ($ws in $wb.Worksheets) {
    $wasProtected = $ws.ProtectContents
    if ($wasProtected) { $ws.Unprotect() }

    # The "!!Compartment" sheet carries an extra top row with the overall
    # "!!!ObjTables ..." workbook-level banner in A1, and its own
    # per-table "!!ObjTables ... id='Compartment' ..." banner in A2.
    # Every other table sheet only has the per-table banner, in A1.
    $a1 = $ws.Range("A1").Text
    if ($a1 -ne $null -and $a1.ToString().Contains("ObjTables")) {
        $ws.Range("A1").Value = Bump-ObjTables($a1.ToString())
    }

    $a2 = $ws.Range("A2").Text
    if ($a2 -ne $null -and $a2.ToString().Contains("ObjTables")) {
        $ws.Range("A2").Value = Bump-ObjTables($a2.ToString())
    }

    if ($wasProtected) { $ws.Protect() }
}

# --- 2) Rename From/To -> FromObject/ToObject -----------------------------

# "!!Relation" sheet: headers in G2 (From) and H2 (To)
$wsRelation = $wb.Worksheets.Item("!!Relation")
$wasProtected = $wsRelation.ProtectContents
if ($wasProtected) { $wsRelation.Unprotect() }

$wsRelation.Range("G2").Value = "!FromObject"
$wsRelation.Range("H2").Value = "!ToObject"

$valG = $wsRelation.Range("G2").Validation
$valG.ErrorTitle = "FromObject"
$valG.InputTitle = "FromObject"

$valH = $wsRelation.Range("H2").Validation
$valH.ErrorTitle = "ToObject"
$valH.InputTitle = "ToObject"

if ($wasProtected) { $wsRelation.Protect() }

# "!!Relationship" sheet: headers in B2 (From) and C2 (To)
$wsRelationship = $wb.Worksheets.Item("!!Relationship")
$wasProtected = $wsRelationship.ProtectContents
if ($wasProtected) { $wsRelationship.Unprotect() }

$wsRelationship.Range("B2").Value = "!FromObject"
$wsRelationship.Range("C2").Value = "!ToObject"

$valB = $wsRelationship.Range("B2").Validation
$valB.ErrorTitle = "FromObject"
$valB.InputTitle = "FromObject"

$valC = $wsRelationship.Range("C2").Validation
$valC.ErrorTitle = "ToObject"
$valC.InputTitle = "ToObject"

if ($wasProtected) { $wsRelationship.Protect() }
